$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (per diff) ---
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(9).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(20).ColumnWidth = 8.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17

# --- Update data rows 2-5 with new values ---
# Row 2
$ws.Cells.Item(2, 1).Value = 45063.50694444445
$ws.Cells.Item(2, 2).Value = 23.541
$ws.Cells.Item(2, 3).Value = 16.304
$ws.Cells.Item(2, 4).Value = 4.244
$ws.Cells.Item(2, 5).Value = 49.627
$ws.Cells.Item(2, 6).Value = 41.058
$ws.Cells.Item(2, 7).Value = 18.526
$ws.Cells.Item(2, 8).Value = 61.718
$ws.Cells.Item(2, 9).Value = 28.505
$ws.Cells.Item(2, 10).Value = 12.158
$ws.Cells.Item(2, 11).Value = 18.744
$ws.Cells.Item(2, 12).Value = 19.59
$ws.Cells.Item(2, 13).Value = 20.45
$ws.Cells.Item(2, 14).Value = 5.915
$ws.Cells.Item(2, 15).Value = 18.422
$ws.Cells.Item(2, 16).Value = 25.954
$ws.Cells.Item(2, 17).Value = 15.331
$ws.Cells.Item(2, 18).Value = 3.814
$ws.Cells.Item(2, 19).Value = 2.604
$ws.Cells.Item(2, 20).Value = 273.066
$ws.Cells.Item(2, 21).Value = 51.313
$ws.Cells.Item(2, 22).Value = 17.004
$ws.Cells.Item(2, 23).Value = 34.141
$ws.Cells.Item(2, 24).Value = 17.76
$ws.Cells.Item(2, 25).Value = 2.295
$ws.Cells.Item(2, 26).Value = 30.618
$ws.Cells.Item(2, 27).Value = 15.02
$ws.Cells.Item(2, 28).Value = 13.452
$ws.Cells.Item(2, 29).Value = 15.735
$ws.Cells.Item(2, 30).Value = 20.304
$ws.Cells.Item(2, 31).Value = 3.64
$ws.Cells.Item(2, 32).Value = 54.593
$ws.Cells.Item(2, 33).Value = 9.503
$ws.Cells.Item(2, 34).Value = 21.259

# Row 3
$ws.Cells.Item(3, 1).Value = 45063.51388888889
$ws.Cells.Item(3, 2).Value = 10.569
$ws.Cells.Item(3, 3).Value = 7.263
$ws.Cells.Item(3, 4).Value = 1.62
$ws.Cells.Item(3, 5).Value = 22.369
$ws.Cells.Item(3, 6).Value = 18.538
$ws.Cells.Item(3, 7).Value = 8.318
$ws.Cells.Item(3, 8).Value = 35.443
$ws.Cells.Item(3, 9).Value = 12.798
$ws.Cells.Item(3, 10).Value = 5.449
$ws.Cells.Item(3, 11).Value = 8.32
$ws.Cells.Item(3, 12).Value = 8.939
$ws.Cells.Item(3, 13).Value = 9.281
$ws.Cells.Item(3, 14).Value = 2.66
$ws.Cells.Item(3, 15).Value = 8.271
$ws.Cells.Item(3, 16).Value = 11.631
$ws.Cells.Item(3, 17).Value = 7.137
$ws.Cells.Item(3, 18).Value = 1.584
$ws.Cells.Item(3, 19).Value = 0.898
$ws.Cells.Item(3, 20).Value = 118.591
$ws.Cells.Item(3, 21).Value = 23.271
$ws.Cells.Item(3, 22).Value = 7.635
$ws.Cells.Item(3, 23).Value = 15.341
$ws.Cells.Item(3, 24).Value = 8.216
$ws.Cells.Item(3, 25).Value = 0.983
$ws.Cells.Item(3, 26).Value = 16.548
$ws.Cells.Item(3, 27).Value = 6.744
$ws.Cells.Item(3, 28).Value = 6.172
$ws.Cells.Item(3, 29).Value = 7.207
$ws.Cells.Item(3, 30).Value = 9.292
$ws.Cells.Item(3, 31).Value = 1.294
$ws.Cells.Item(3, 32).Value = 32.113
$ws.Cells.Item(3, 33).Value = 4.203
$ws.Cells.Item(3, 34).Value = 9.546

# Row 4
$ws.Cells.Item(4, 1).Value = 45063.52083333334
$ws.Cells.Item(4, 2).Value = 5.765
$ws.Cells.Item(4, 3).Value = 3.883
$ws.Cells.Item(4, 4).Value = 0.967
$ws.Cells.Item(4, 5).Value = 12.135
$ws.Cells.Item(4, 6).Value = 10.071
$ws.Cells.Item(4, 7).Value = 4.537
$ws.Cells.Item(4, 8).Value = 20.531
$ws.Cells.Item(4, 9).Value = 6.981
$ws.Cells.Item(4, 10).Value = 2.943
$ws.Cells.Item(4, 11).Value = 4.475
$ws.Cells.Item(4, 12).Value = 4.898
$ws.Cells.Item(4, 13).Value = 5.019
$ws.Cells.Item(4, 14).Value = 1.453
$ws.Cells.Item(4, 15).Value = 4.512
$ws.Cells.Item(4, 16).Value = 6.316
$ws.Cells.Item(4, 17).Value = 3.996
$ws.Cells.Item(4, 18).Value = 0.998
$ws.Cells.Item(4, 19).Value = 0.499
$ws.Cells.Item(4, 20).Value = 61.358
$ws.Cells.Item(4, 21).Value = 12.742
$ws.Cells.Item(4, 22).Value = 4.164
$ws.Cells.Item(4, 23).Value = 8.325
$ws.Cells.Item(4, 24).Value = 4.551
$ws.Cells.Item(4, 25).Value = 0.503
$ws.Cells.Item(4, 26).Value = 9.351
$ws.Cells.Item(4, 27).Value = 3.678
$ws.Cells.Item(4, 28).Value = 3.417
$ws.Cells.Item(4, 29).Value = 3.982
$ws.Cells.Item(4, 30).Value = 5.082
$ws.Cells.Item(4, 31).Value = 0.784
$ws.Cells.Item(4, 32).Value = 18.546
$ws.Cells.Item(4, 33).Value = 2.252
$ws.Cells.Item(4, 34).Value = 5.207

# Row 5
$ws.Cells.Item(5, 1).Value = 45063.52777777778
$ws.Cells.Item(5, 2).Value = 20.66
$ws.Cells.Item(5, 3).Value = 15.19
$ws.Cells.Item(5, 4).Value = 1.25
$ws.Cells.Item(5, 5).Value = 44.63
$ws.Cells.Item(5, 6).Value = 36.94
$ws.Cells.Item(5, 7).Value = 16.26
$ws.Cells.Item(5, 8).Value = 59.47
$ws.Cells.Item(5, 9).Value = 25.01
$ws.Cells.Item(5, 10).Value = 11.07
$ws.Cells.Item(5, 11).Value = 16.61
$ws.Cells.Item(5, 12).Value = 17.96
$ws.Cells.Item(5, 13).Value = 18.85
$ws.Cells.Item(5, 14).Value = 5.19
$ws.Cells.Item(5, 15).Value = 16.17
$ws.Cells.Item(5, 16).Value = 22.98
$ws.Cells.Item(5, 17).Value = 13.61
$ws.Cells.Item(5, 18).Value = 0.87
$ws.Cells.Item(5, 19).Value = 0.84
$ws.Cells.Item(5, 20).Value = 238.74
$ws.Cells.Item(5, 21).Value = 45.12
$ws.Cells.Item(5, 22).Value = 14.92
$ws.Cells.Item(5, 23).Value = 30.31
$ws.Cells.Item(5, 24).Value = 16.13
$ws.Cells.Item(5, 25).Value = 2.04
$ws.Cells.Item(5, 26).Value = 29.3
$ws.Cells.Item(5, 27).Value = 13.18
$ws.Cells.Item(5, 28).Value = 11.71
$ws.Cells.Item(5, 29).Value = 13.75
$ws.Cells.Item(5, 30).Value = 18.83
$ws.Cells.Item(5, 31).Value = 0.56
$ws.Cells.Item(5, 32).Value = 53.64
$ws.Cells.Item(5, 33).Value = 8.41
$ws.Cells.Item(5, 34).Value = 18.66

# --- Remove row 6 (data row no longer present) ---
$ws.Rows.Item(6).Delete()

Write-Host "Edit complete"
